$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last 4 rows (rows 14-17), which correspond to the
# "Resolving-Mac" sending-cluster block that no longer exists in the
# recomputed TPM data (dimension shrinks from A1:T17 to A1:T13).
$ws.Range("A14:T17").EntireRow.Delete() | Out-Null

# Update the remaining 12 data rows (rows 2-13) with the recomputed
# TPM-derived values, including the shuffled "Target cluster" pairing
# (MuSCs is no longer present as a target cluster).
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Bgn"
$ws.Cells.Item(2, 3).Value = "Tlr2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 99.28451533333333
$ws.Cells.Item(2, 8).Value = 297.853546
$ws.Cells.Item(2, 9).Value = 0.02270354261926982
$ws.Cells.Item(2, 10).Value = 0.02270354261926982
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 1.597802666666666
$ws.Cells.Item(2, 14).Value = 4.793407999999999
$ws.Cells.Item(2, 15).Value = 0.02304920886321625
$ws.Cells.Item(2, 16).Value = 0.02304920886321625
$ws.Cells.Item(2, 17).Value = 158.6370633583075
$ws.Cells.Item(2, 18).Value = 1427.733570224768
$ws.Cells.Item(2, 19).Value = 0.0005232986957664818
$ws.Cells.Item(2, 20).Value = 0.0005232986957664818
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Bgn"
$ws.Cells.Item(3, 3).Value = "Tlr2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 99.28451533333333
$ws.Cells.Item(3, 8).Value = 297.853546
$ws.Cells.Item(3, 9).Value = 0.02270354261926982
$ws.Cells.Item(3, 10).Value = 0.02270354261926982
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 1.27306
$ws.Cells.Item(3, 14).Value = 3.81918
$ws.Cells.Item(3, 15).Value = 0.0183646118807784
$ws.Cells.Item(3, 16).Value = 0.0183646118807784
$ws.Cells.Item(3, 17).Value = 126.3951450902533
$ws.Cells.Item(3, 18).Value = 1137.55630581228
$ws.Cells.Item(3, 19).Value = 0.0004169417485216014
$ws.Cells.Item(3, 20).Value = 0.0004169417485216014
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Bgn"
$ws.Cells.Item(4, 3).Value = "Tlr2"
$ws.Cells.Item(4, 4).Value = "Resolving-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 99.28451533333333
$ws.Cells.Item(4, 8).Value = 297.853546
$ws.Cells.Item(4, 9).Value = 0.02270354261926982
$ws.Cells.Item(4, 10).Value = 0.02270354261926982
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 66.45050433333334
$ws.Cells.Item(4, 14).Value = 199.351513
$ws.Cells.Item(4, 15).Value = 0.9585861792560053
$ws.Cells.Item(4, 16).Value = 0.9585861792560054
$ws.Cells.Item(4, 17).Value = 6597.506116390568
$ws.Cells.Item(4, 18).Value = 59377.5550475151
$ws.Cells.Item(4, 19).Value = 0.02176330217498174
$ws.Cells.Item(4, 20).Value = 0.02176330217498174
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Bgn"
$ws.Cells.Item(5, 3).Value = "Tlr2"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 4010.868571
$ws.Cells.Item(5, 8).Value = 12032.605713
$ws.Cells.Item(5, 9).Value = 0.9171714767027319
$ws.Cells.Item(5, 10).Value = 0.9171714767027318
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 1.597802666666666
$ws.Cells.Item(5, 14).Value = 4.793407999999999
$ws.Cells.Item(5, 15).Value = 0.02304920886321625
$ws.Cells.Item(5, 16).Value = 0.02304920886321625
$ws.Cells.Item(5, 17).Value = 6408.576498393322
$ws.Cells.Item(5, 18).Value = 57677.1884855399
$ws.Cells.Item(5, 19).Value = 0.02114007692990574
$ws.Cells.Item(5, 20).Value = 0.02114007692990574
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Bgn"
$ws.Cells.Item(6, 3).Value = "Tlr2"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 4010.868571
$ws.Cells.Item(6, 8).Value = 12032.605713
$ws.Cells.Item(6, 9).Value = 0.9171714767027319
$ws.Cells.Item(6, 10).Value = 0.9171714767027318
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 1.27306
$ws.Cells.Item(6, 14).Value = 3.81918
$ws.Cells.Item(6, 15).Value = 0.0183646118807784
$ws.Cells.Item(6, 16).Value = 0.0183646118807784
$ws.Cells.Item(6, 17).Value = 5106.07634299726
$ws.Cells.Item(6, 18).Value = 45954.68708697534
$ws.Cells.Item(6, 19).Value = 0.01684349819776606
$ws.Cells.Item(6, 20).Value = 0.01684349819776606
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Bgn"
$ws.Cells.Item(7, 3).Value = "Tlr2"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 4010.868571
$ws.Cells.Item(7, 8).Value = 12032.605713
$ws.Cells.Item(7, 9).Value = 0.9171714767027319
$ws.Cells.Item(7, 10).Value = 0.9171714767027318
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 66.45050433333334
$ws.Cells.Item(7, 14).Value = 199.351513
$ws.Cells.Item(7, 15).Value = 0.9585861792560053
$ws.Cells.Item(7, 16).Value = 0.9585861792560054
$ws.Cells.Item(7, 17).Value = 266524.2393576661
$ws.Cells.Item(7, 18).Value = 2398718.154218994
$ws.Cells.Item(7, 19).Value = 0.8791879015750601
$ws.Cells.Item(7, 20).Value = 0.8791879015750601
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Bgn"
$ws.Cells.Item(8, 3).Value = "Tlr2"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 249.2612966666667
$ws.Cells.Item(8, 8).Value = 747.7838899999999
$ws.Cells.Item(8, 9).Value = 0.05699896356653876
$ws.Cells.Item(8, 10).Value = 0.05699896356653875
$ws.Cells.Item(8, 11).Value = 2
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 1.597802666666666
$ws.Cells.Item(8, 14).Value = 4.793407999999999
$ws.Cells.Item(8, 15).Value = 0.02304920886321625
$ws.Cells.Item(8, 16).Value = 0.02304920886321625
$ws.Cells.Item(8, 17).Value = 398.2703645107911
$ws.Cells.Item(8, 18).Value = 3584.433280597119
$ws.Cells.Item(8, 19).Value = 0.001313781016232005
$ws.Cells.Item(8, 20).Value = 0.001313781016232005
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Bgn"
$ws.Cells.Item(9, 3).Value = "Tlr2"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 249.2612966666667
$ws.Cells.Item(9, 8).Value = 747.7838899999999
$ws.Cells.Item(9, 9).Value = 0.05699896356653876
$ws.Cells.Item(9, 10).Value = 0.05699896356653875
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 1.27306
$ws.Cells.Item(9, 14).Value = 3.81918
$ws.Cells.Item(9, 15).Value = 0.0183646118807784
$ws.Cells.Item(9, 16).Value = 0.0183646118807784
$ws.Cells.Item(9, 17).Value = 317.3245863344666
$ws.Cells.Item(9, 18).Value = 2855.921277010199
$ws.Cells.Item(9, 19).Value = 0.001046763843506113
$ws.Cells.Item(9, 20).Value = 0.001046763843506113
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Bgn"
$ws.Cells.Item(10, 3).Value = "Tlr2"
$ws.Cells.Item(10, 4).Value = "Resolving-Mac"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 249.2612966666667
$ws.Cells.Item(10, 8).Value = 747.7838899999999
$ws.Cells.Item(10, 9).Value = 0.05699896356653876
$ws.Cells.Item(10, 10).Value = 0.05699896356653875
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 66.45050433333334
$ws.Cells.Item(10, 14).Value = 199.351513
$ws.Cells.Item(10, 15).Value = 0.9585861792560053
$ws.Cells.Item(10, 16).Value = 0.9585861792560054
$ws.Cells.Item(10, 17).Value = 16563.53887428062
$ws.Cells.Item(10, 18).Value = 149071.8498685256
$ws.Cells.Item(10, 19).Value = 0.05463841870680064
$ws.Cells.Item(10, 20).Value = 0.05463841870680064
$ws.Cells.Item(11, 1).Value = "Resolving-Mac"
$ws.Cells.Item(11, 2).Value = "Bgn"
$ws.Cells.Item(11, 3).Value = "Tlr2"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 13.67033766666667
$ws.Cells.Item(11, 8).Value = 41.011013
$ws.Cells.Item(11, 9).Value = 0.003126017111459632
$ws.Cells.Item(11, 10).Value = 0.003126017111459632
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 12).Value = 0.6666666666666666
$ws.Cells.Item(11, 13).Value = 1.597802666666666
$ws.Cells.Item(11, 14).Value = 4.793407999999999
$ws.Cells.Item(11, 15).Value = 0.02304920886321625
$ws.Cells.Item(11, 16).Value = 0.02304920886321625
$ws.Cells.Item(11, 17).Value = 21.84250197803378
$ws.Cells.Item(11, 18).Value = 196.582517802304
$ws.Cells.Item(11, 19).Value = 0.000072052221312021
$ws.Cells.Item(11, 20).Value = 0.000072052221312021
$ws.Cells.Item(12, 1).Value = "Resolving-Mac"
$ws.Cells.Item(12, 2).Value = "Bgn"
$ws.Cells.Item(12, 3).Value = "Tlr2"
$ws.Cells.Item(12, 4).Value = "FAPs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 13.67033766666667
$ws.Cells.Item(12, 8).Value = 41.011013
$ws.Cells.Item(12, 9).Value = 0.003126017111459632
$ws.Cells.Item(12, 10).Value = 0.003126017111459632
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 1.27306
$ws.Cells.Item(12, 14).Value = 3.81918
$ws.Cells.Item(12, 15).Value = 0.0183646118807784
$ws.Cells.Item(12, 16).Value = 0.0183646118807784
$ws.Cells.Item(12, 17).Value = 17.40316006992666
$ws.Cells.Item(12, 18).Value = 156.62844062934
$ws.Cells.Item(12, 19).Value = 0.00005740809098462815
$ws.Cells.Item(12, 20).Value = 0.00005740809098462815
$ws.Cells.Item(13, 1).Value = "Resolving-Mac"
$ws.Cells.Item(13, 2).Value = "Bgn"
$ws.Cells.Item(13, 3).Value = "Tlr2"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 13.67033766666667
$ws.Cells.Item(13, 8).Value = 41.011013
$ws.Cells.Item(13, 9).Value = 0.003126017111459632
$ws.Cells.Item(13, 10).Value = 0.003126017111459632
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 66.45050433333334
$ws.Cells.Item(13, 14).Value = 199.351513
$ws.Cells.Item(13, 15).Value = 0.9585861792560053
$ws.Cells.Item(13, 16).Value = 0.9585861792560054
$ws.Cells.Item(13, 17).Value = 908.4008323569633
$ws.Cells.Item(13, 18).Value = 8175.607491212669
$ws.Cells.Item(13, 19).Value = 0.002996556799162983
$ws.Cells.Item(13, 20).Value = 0.002996556799162983
